$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Drop every existing hyperlink object (the runtime's Hyperlinks.Delete()
#    clears the whole worksheet collection regardless of which Range it was
#    scoped from, so a single call on any range is enough). We rebuild the
#    hyperlinks for the surviving rows further down.
# ---------------------------------------------------------------------------
$ws.Range("A1").Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 2. Remove the old rows 12-15 entirely (the refreshed scrape only keeps 10
#    data rows now, i.e. rows 2-11).
# ---------------------------------------------------------------------------
$ws.Range("A12:H15").Clear()

# ---------------------------------------------------------------------------
# 3. Overwrite rows 2-11 with the newly scraped listings.
# ---------------------------------------------------------------------------
$rows = @(
    @{ Row=2;  A="2025-10-14 06:28:19"; B="Amazon商品を自動抽出してBASEに出品するツール開発(スクレイピング機能)"; C="システム開発"; D="20,000 円 ~ 50,000 円 / 固定";        E="期限情報なし"; F="https://www.lancers.jp/work/detail/5412467"; G=168; H="◆ツール,開発" },
    @{ Row=3;  A="2025-10-14 06:28:19"; B="海外仕入れ元サイト→ツールを動かす為のCSVファイルに週1で自動抽出の制作(自動/スクレイピング)"; C="システム開発"; D="5,000 円 ~ 10,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5251319"; G=135; H="◆ツール,スクレイピング ◇サイト" },
    @{ Row=4;  A="2025-10-14 06:28:19"; B="【フルリモート】WordPress開発スタッフ募集";                              C="システム開発"; D="50,000 円 ~ 100,000 円 / 固定";       E="期限情報なし"; F="https://www.lancers.jp/work/detail/5407811"; G=88;  H="◆開発 ○WordPress" },
    @{ Row=5;  A="2025-10-14 06:28:19"; B="FileMaker開発";                                                       C="システム開発"; D="100,000 円 ~ 200,000 円 / 固定";      E="期限情報なし"; F="https://www.lancers.jp/work/detail/5412487"; G=68;  H="◆開発" },
    @{ Row=6;  A="2025-10-14 06:28:19"; B="初回 ポケパラの自動いいね等の開発";                                      C="システム開発"; D="20,000 円 ~ 50,000 円 / 固定";       E="期限情報なし"; F="https://www.lancers.jp/work/detail/5412453"; G=63;  H="◆開発" },
    @{ Row=7;  A="2025-10-14 06:28:19"; B="【急募】Salesforce・MA・CRMコンサルタント経験者を探しています!";          C="システム開発"; D="200,000 円 ~ 300,000 円 / 固定";      E="期限情報なし"; F="https://www.lancers.jp/work/detail/5371747"; G=48;  H="◆コンサル" },
    @{ Row=8;  A="2025-10-14 06:28:19"; B="【急募】モバイルアプリ テスト業務 委託募集(3 - 4週間)";                  C="システム開発"; D="100,000 円 ~ 200,000 円 / 固定";      E="期限情報なし"; F="https://www.lancers.jp/work/detail/5412563"; G=38;  H="◇アプリ" },
    @{ Row=9;  A="2025-10-14 06:28:19"; B="【急募】警備スタッフマッチングシステム構築の依頼";                        C="システム開発"; D="1,000,000 円 ~ 3,000,000 円 / 固定";  E="期限情報なし"; F="https://www.lancers.jp/work/detail/5412802"; G=40;  H=$null },
    @{ Row=10; A="2025-10-14 06:28:19"; B="Drupal関連プロジェクトの要件定義や基本設計ができる方(1人月、長期継続案件)"; C="システム開発"; D="300,000 円 ~ 500,000 円 / 固定";       E="期限情報なし"; F="https://www.lancers.jp/work/detail/5400683"; G=25;  H=$null },
    @{ Row=11; A="2025-10-14 06:28:19"; B="【急募】【高単価】赤坂で15分のコンテンツ更新作業をお手伝いください!";      C="システム開発"; D="20,000 円 ~ 50,000 円 / 固定";       E="期限情報なし"; F="https://www.lancers.jp/work/detail/5412531"; G=13;  H=$null }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
    if ($r.H -eq $null) {
        $ws.Cells.Item($r.Row, 8).ClearContents()
    } else {
        $ws.Cells.Item($r.Row, 8).Value = $r.H
    }
}

# ---------------------------------------------------------------------------
# 4. Re-create the hyperlinks on column F for the surviving rows 2-11,
#    pointing at the refreshed URLs (and restore the plain "Hyperlink" cell
#    style so it matches s="1" like the rest of the workbook).
# ---------------------------------------------------------------------------
foreach ($r in $rows) {
    $ws.Hyperlinks.Add($ws.Range("F" + $r.Row), $r.F)
    $ws.Range("F" + $r.Row).Style = "Hyperlink"
}

# ---------------------------------------------------------------------------
# 5. Column width tweaks (stored width = input + 0.8333... in this engine,
#    so back the input off to land exactly on 32 / 19).
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 31.166666666666668
$ws.Columns.Item(8).ColumnWidth = 18.166666666666668
